$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "User 2 added something in cell A4"
$ws.Range("A5").Select()
